$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Chris Paul", "PG", "San Antonio Spurs"),
    @("Keon Johnson", "PG,SG", "Brooklyn Nets"),
    @("Russell Westbrook", "PG,SG", "Denver Nuggets"),
    @("Khris Middleton", "SF", "Milwaukee Bucks"),
    @("Paolo Banchero", "SF,PF", "Orlando Magic"),
    @("Nikola Jokic", "C", "Denver Nuggets"),
    @("Jakob Poeltl", "C", "Toronto Raptors"),
    @("Rudy Gobert", "C", "Minnesota Timberwolves"),
    @("Pascal Siakam", "SF,PF,C", "Indiana Pacers"),
    @("Jaylen Brown", "SG,SF", "Boston Celtics"),
    @("Caris LeVert", "SG,SF", "Cleveland Cavaliers"),
    @("Grayson Allen", "PG,SG,SF", "Phoenix Suns"),
    @("Jalen Green", "PG,SG", "Houston Rockets"),
    @("Buddy Hield", "SG,SF", "Golden State Warriors"),
    @("Dejounte Murray", "PG,SG", "New Orleans Pelicans"),
    @("Chet Holmgren", "PF,C", "Oklahoma City Thunder"),
    @("Jalen Suggs", "PG,SG", "Orlando Magic"),
    @("Deni Avdija", "SF,PF", "Portland Trail Blazers")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
